# Creating a common data provider:
# Add a new "OpenAccountTest" worksheet (after the existing "AddCustomerTest"
# sheet) that holds a small customer/currency lookup table, and widen
# column D on "AddCustomerTest" so the "alerttext" column content fits.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("AddCustomerTest")

# New sheet, placed right after AddCustomerTest.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "OpenAccountTest"

$ws2.Range("A1").Value = "customer"
$ws2.Range("B1").Value = "currency"
$ws2.Range("A2").Value = "Deepender Singh"
$ws2.Range("B2").Value = "Rupee"
$ws2.Range("B2").Select()

# Widen column D on AddCustomerTest to fit its longest value
# ("Customer added successfully").
$ws1.Columns.Item(4).ColumnWidth = 26.3

# Keep the original sheet active/selected (adding a sheet would otherwise
# switch focus to it).
$ws1.Activate()
